$p = $ppt.ActivePresentation

# --- Edit: slide 16 ("ai/6.1" deck) - fix typo "Se dejamos" -> "Si dejamos" ---
$s16 = $p.Slides.Item(16)
$shp369 = $s16.Shapes.Item(2)
$shp369.TextFrame.TextRange.Text = "Si dejamos en blanco el primer o el último número de la rebanada, se asume que es el inicio o el final de la cadena, respectivamente"
